# Add new power plant types to the Electricity Source subscript
# (BDTPTUMCF - Boolean Does This Plant Type Use Maximum Capacity Factor)
# Issues #280 and #99

$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("BDTPTUMCF")

$newPlantTypes = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$startRow = 19
$ws.Cells.Item(18, 2).Copy()
for ($i = 0; $i -lt $newPlantTypes.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newPlantTypes[$i]
    $ws.Cells.Item($row, 2).Value = 1
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
}

$ws.Range("A25").Select()
$wsAbout.Activate()
